# Seed the database: normalize all due_date values in the "bill" sheet to a
# single fixed timestamp, then reset the active selection to D1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All due_date rows (D2:D17) get the same seeded timestamp.
$ws.Range("D2:D17").Value = 45147.012407407405

# Reset selection from F1:F1048576 (whole column F) to just D1.
$ws.Range("D1").Select()
